$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 2067
$ws.Range("J3").Value = 2177
$ws.Range("J4").Value = 492
$ws.Range("I6").Value = 8969
$ws.Range("J6").Value = 2719
$ws.Range("I7").Value = 26201
$ws.Range("J7").Value = 7611

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J7").Value = 221
$ws.Range("J8").Value = 476
$ws.Range("J10").Value = 47
$ws.Range("J11").Value = 105
$ws.Range("J12").Value = 17
$ws.Range("J15").Value = 99
$ws.Range("J19").Value = 252
$ws.Range("J20").Value = 159
$ws.Range("J25").Value = 45
$ws.Range("J29").Value = 426
$ws.Range("J30").Value = 29
$ws.Range("J33").Value = 317
$ws.Range("I35").Value = 34
$ws.Range("J36").Value = 115
$ws.Range("J37").Value = 254
$ws.Range("J41").Value = 46
$ws.Range("J42").Value = 288
$ws.Range("J43").Value = 77
$ws.Range("J44").Value = 61
$ws.Range("J46").Value = 26
$ws.Range("J47").Value = 70
$ws.Range("J48").Value = 70
$ws.Range("J49").Value = 45
$ws.Range("J52").Value = 180
$ws.Range("J53").Value = 73
$ws.Range("J54").Value = 155
$ws.Range("J55").Value = 91
$ws.Range("J63").Value = 28
$ws.Range("J64").Value = 49
$ws.Range("J65").Value = 198
$ws.Range("J66").Value = 16
$ws.Range("J67").Value = 279
$ws.Range("J70").Value = 15
$ws.Range("J73").Value = 69
$ws.Range("J74").Value = 12
$ws.Range("J76").Value = 113
$ws.Range("J77").Value = 55
$ws.Range("J78").Value = 97
$ws.Range("J79").Value = 227
$ws.Range("J83").Value = 183
$ws.Range("J84").Value = 76
$ws.Range("J85").Value = 357
$ws.Range("J92").Value = 25
$ws.Range("J93").Value = 34
$ws.Range("J94").Value = 58
$ws.Range("J95").Value = 115
$ws.Range("J96").Value = 84
$ws.Range("J98").Value = 49
$ws.Range("J99").Value = 103
$ws.Range("J100").Value = 16
$ws.Range("I101").Value = 26201
$ws.Range("J101").Value = 7611

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("I6").Value = 14
$ws.Range("I7").Value = 34

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J3").Value = 23
$ws.Range("J6").Value = 18
$ws.Range("J7").Value = 69

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J3").Value = 24
$ws.Range("J7").Value = 84

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("J3").Value = 9
$ws.Range("J7").Value = 29

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 73
$ws.Range("J3").Value = 91
$ws.Range("J7").Value = 254

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J3").Value = 35
$ws.Range("J4").Value = 7
$ws.Range("J7").Value = 103

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 57
$ws.Range("J3").Value = 109
$ws.Range("J7").Value = 279

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J6").Value = 29
$ws.Range("J7").Value = 76

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 59
$ws.Range("J7").Value = 198

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J3").Value = 62
$ws.Range("J6").Value = 53
$ws.Range("J7").Value = 183

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J2").Value = 41
$ws.Range("J7").Value = 115

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 80
$ws.Range("J4").Value = 15
$ws.Range("J6").Value = 112
$ws.Range("J7").Value = 317

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J2").Value = 9
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 45

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J3").Value = 28
$ws.Range("J7").Value = 155

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 128
$ws.Range("J3").Value = 148
$ws.Range("J6").Value = 113
$ws.Range("J7").Value = 426

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 63
$ws.Range("J6").Value = 94
$ws.Range("J7").Value = 252

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J3").Value = 17
$ws.Range("J6").Value = 18
$ws.Range("J7").Value = 61

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J6").Value = 35
$ws.Range("J7").Value = 70

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J3").Value = 26
$ws.Range("J6").Value = 64
$ws.Range("J7").Value = 113

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J3").Value = 144
$ws.Range("J6").Value = 99
$ws.Range("J7").Value = 357

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J6").Value = 21
$ws.Range("J7").Value = 46

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J3").Value = 61
$ws.Range("J6").Value = 145
$ws.Range("J7").Value = 288

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("J2").Value = 15
$ws.Range("J7").Value = 47

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J2").Value = 24
$ws.Range("J7").Value = 97

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J3").Value = 15
$ws.Range("J7").Value = 91

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("J2").Value = 7
$ws.Range("J7").Value = 26

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J3").Value = 85
$ws.Range("J6").Value = 61
$ws.Range("J7").Value = 227

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("J3").Value = 12
$ws.Range("J7").Value = 49

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J3").Value = 50
$ws.Range("J6").Value = 43
$ws.Range("J7").Value = 159

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J3").Value = 28
$ws.Range("J6").Value = 45
$ws.Range("J7").Value = 115

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 34

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("J5").Value = 7
$ws.Range("J6").Value = 16

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J3").Value = 56
$ws.Range("J7").Value = 180

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J6").Value = 34
$ws.Range("J7").Value = 58

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J6").Value = 8
$ws.Range("J7").Value = 45

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J3").Value = 20
$ws.Range("J6").Value = 31
$ws.Range("J7").Value = 70

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J3").Value = 24
$ws.Range("J6").Value = 46
$ws.Range("J7").Value = 99

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("J6").Value = 27
$ws.Range("J7").Value = 49

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("J2").Value = 2
$ws.Range("J7").Value = 16

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J6").Value = 39
$ws.Range("J7").Value = 105

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 25

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("J3").Value = 4
$ws.Range("J7").Value = 15

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 150
$ws.Range("J3").Value = 156
$ws.Range("J6").Value = 136
$ws.Range("J7").Value = 476

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J6").Value = 47
$ws.Range("J7").Value = 77

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J6").Value = 43
$ws.Range("J7").Value = 73

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("J4").Value = 6
$ws.Range("J7").Value = 55

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 71
$ws.Range("J3").Value = 69
$ws.Range("J7").Value = 221

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("J6").Value = 13
$ws.Range("J7").Value = 17

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("J3").Value = 6
$ws.Range("J7").Value = 12
